$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header values for columns G, H, I (row 2)
$ws.Range("G2").Value = "WIN_IP_Address"
$ws.Range("H2").Value = "USERNAME"
$ws.Range("I2").Value = "PASSWORD"

# New body rows (rows 3-7): WIN IP address / username / password
$ws.Range("G3").Value = "13.201.47.117"
$ws.Range("H3").Value = "Administrator"
$ws.Range("I3").Value = "g)HQzIo5pD*mdyXb.p6NxwerZ-EpLdz?"

$ws.Range("G4").Value = "13.232.8.245"
$ws.Range("H4").Value = "Administrator"
$ws.Range("I4").Value = "b4hv?p`$Jgi!-PvMdx?K7?y!;pf=aH&yb"

$ws.Range("G5").Value = "52.66.201.47"
$ws.Range("H5").Value = "Administrator"
$ws.Range("I5").Value = "&a`$Qr)OYMb2uR0N@tiTpusiEVszgPesd"

$ws.Range("G6").Value = "65.0.205.20"
$ws.Range("H6").Value = "Administrator"
$ws.Range("I6").Value = "HJ=*.LybO?y6fV6)kve)2o%TDB?D7?.z"

$ws.Range("G7").Value = "13.127.48.86"
$ws.Range("H7").Value = "Administrator"
$ws.Range("I7").Value = "hi.-pqU)8.eJFoMxAr`$ij7lY2lg)GrAi"

# Give the new header cells (G2:I2) the same bold/wrap/vcenter look as the
# rest of row 2, by copying the existing header formatting over.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("G2:I2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column widths for the new columns
$ws.Range("G1").ColumnWidth = 25.8333
$ws.Range("H1").ColumnWidth = 19.666
$ws.Range("I1").ColumnWidth = 50.666

# Selection ends on O6
$ws.Range("O6").Select()
